$wb = $excel.ActiveWorkbook

# --- Accounts sheet: rename accounts to include the account number suffix ---
$wsAccounts = $wb.Worksheets.Item("Accounts")
$wsAccounts.Range("B6").Value = "Credit 7891"
$wsAccounts.Range("B5").Value = "Mortgage 7890"
$wsAccounts.Range("B2").Value = "Checking 1111"
$wsAccounts.Range("B3").Value = "Savings 2222"
$wsAccounts.Range("B4").Value = "Loan 3333"
$wsAccounts.Columns.Item(2).AutoFit() | Out-Null

# --- AccountsInfo sheet: same renaming + eligibility flags for the Mortgage row ---
$wsAccountsInfo = $wb.Worksheets.Item("AccountsInfo")
$wsAccountsInfo.Range("C2").Value = "Checking 1111"
$wsAccountsInfo.Range("C3").Value = "Mortgage 7890"
$wsAccountsInfo.Range("D3").Value = $false
$wsAccountsInfo.Range("E3").Value = $false
$wsAccountsInfo.Range("C4").Value = "Credit 7891"

# --- Selections / active-tab bookkeeping ---
$wsAccounts.Range("B4").Select()
$wsAccountsInfo.Range("D3").Select()
